# Update "想去人数" (want-to-go count, column F) figures to the freshly
# scraped values as published to gh-pages at commit 456a3b4.
#
# Sheet1 "展览" (Exhibition), Sheet2 "演出" (Performance) and
# Sheet4 "全部类型" (All types) each contain F-column values that need to be
# refreshed; Sheet4 simply aggregates the rows from the other sheets, which
# is why most values repeat across sheets.

function Set-ColF {
    param($Sheet, $RowValues)
    foreach ($row in $RowValues.Keys) {
        $Sheet.Cells.Item($row, 6).Value = $RowValues[$row]
    }
}

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) -> sheet1
$ws1 = $wb.Worksheets.Item("展览")
$rows1 = @{
    3  = 1160
    4  = 12568
    5  = 696
    10 = 316
    11 = 1841
    15 = 204
    17 = 326
    18 = 221
    20 = 111
    21 = 117
    22 = 25
    23 = 205
    24 = 234
    25 = 1249
    26 = 58
}
Set-ColF $ws1 $rows1

# Sheet "演出" (Performance) -> sheet2
$ws2 = $wb.Worksheets.Item("演出")
$rows2 = @{
    4  = 271
    6  = 126
    10 = 345
}
Set-ColF $ws2 $rows2

# Sheet "全部类型" (All types) -> sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$rows4 = @{
    2  = 858
    6  = 1160
    7  = 12568
    8  = 271
    9  = 696
    14 = 316
    15 = 1841
    20 = 204
    21 = 126
    22 = 126
    27 = 345
    28 = 326
    30 = 221
    32 = 111
    33 = 117
    34 = 25
    36 = 205
    39 = 234
    40 = 1249
    42 = 58
}
Set-ColF $ws4 $rows4
